$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 4.715800000000001
$ws.Range("B10").Value = 8.509300000000005
$ws.Range("B12").Value = 6.313400000000001
$ws.Range("E13").Value = 12.1886
$ws.Range("B18").Value = 5.191500000000006
